$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = -1.1684498070942733
$ws.Range("C3").Value = -1.0170928201827825

$ws.Range("B4").Value = -1.2288377017848569
$ws.Range("C4").Value = -1.0308273760089839

$ws.Range("B5").Value = -1.2268230227574015
$ws.Range("C5").Value = -0.99852773427268726

$ws.Range("B6").Value = -1.3436354936910277
$ws.Range("C6").Value = -1.0900851267412499

$ws.Range("B7").Value = -1.3176761516052331
$ws.Range("C7").Value = -1.0442165601249487

$ws.Range("B9").Value = 0.12261107938191873
$ws.Range("C9").Value = 0.11417102325186967

$ws.Range("B13").Value = 0.061230399542213307
$ws.Range("C13").Value = 0.085500224730436808

$ws.Range("B14").Value = 1.727188947827927
$ws.Range("C14").Value = 1.7554254256659443

$ws.Range("C16").Value = -0.2427844851356811
$ws.Range("C17").Value = 0.92754283549268335
$ws.Range("C18").Value = 0.29709434521137013
$ws.Range("C19").Value = -1.5533176404603715
$ws.Range("C20").Value = 0.0042060602902282843

$ws.Range("B21").Value = 8415
$ws.Range("C21").Value = 10251
